$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting from an existing header cell (AC1) to the new headers
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in team record (Wins/Losses/Ties) for every data row
$ws.Range("AD2:AD51").Value = 72
$ws.Range("AE2:AE51").Value = 90
$ws.Range("AF2:AF51").Value = 0
